# Add newly uploaded webcam location rows to the "location-1" worksheet.
# Column layout (row 1 headers): A=Category, B=latitude/longitude,
# C=Location, D=CITY, E=COUNTRY, F=YouTube Link, G=Status(formula, not used
# for these new rows - matches rows 50-53 above them).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 54 - Kaohsiung, Taiwan
$ws.Cells.Item(54, 6).Value = "C03Itx8iSC0"
$ws.Cells.Item(54, 5).Value = "Taiwan"
$ws.Cells.Item(54, 4).Value = "Kaohsiung"
$ws.Cells.Item(54, 3).Value = "高雄壽山情人觀景台 4K即時影像 | Kaohsiung Shoushan Lovers' Viewing Platform 4K Live Camera"
$ws.Cells.Item(54, 2).Value = "22.625447377485727, 120.27414133578823"
$ws.Cells.Item(54, 1).Value = "LIVE, SEA, PORT, SHIP"

# Row 55 - Taipei, Taiwan
$ws.Cells.Item(55, 6).Value = "z_fY1pj1VBw"
$ws.Cells.Item(55, 4).Value = "Taipei"
$ws.Cells.Item(55, 2).Value = "25.033763020202137, 121.56503162459309"
$ws.Cells.Item(55, 3).Value = "Taipei City panorama view - TAIPEI 101"
$ws.Cells.Item(55, 1).Value = "LIVE, CITY, BUILDING"
$ws.Cells.Item(55, 5).Value = "Taiwan"

# Row 56 - Zaanse Schans, Netherlands
$ws.Cells.Item(56, 6).Value = "o9MIV7sep5k"
$ws.Cells.Item(56, 4).Value = "Zaanse Schans"
$ws.Cells.Item(56, 5).Value = "Netherlands"
$ws.Cells.Item(56, 3).Value = "WebCam.NL | dezaanseschans.nl | live ultraHD Pan Tilt Zoom camera"
$ws.Cells.Item(56, 2).Value = "52.475571753497086, 4.818457427726254"
$ws.Cells.Item(56, 1).Value = "LIVE, LANDSCAPE, RIVER"

# Match the bordered style used by columns A, C, D, E in the rows directly
# above (B and F keep the default/no style, same as rows 50-53).
for ($r = 54; $r -le 56; $r++) {
    $ws.Cells.Item($r, 1).Style = $ws.Cells.Item($r - 1, 1).Style
    $ws.Cells.Item($r, 3).Style = $ws.Cells.Item($r - 1, 3).Style
    $ws.Cells.Item($r, 4).Style = $ws.Cells.Item($r - 1, 4).Style
    $ws.Cells.Item($r, 5).Style = $ws.Cells.Item($r - 1, 5).Style
}

# Reflect the author's post-upload scroll/selection position.
$excel.Goto($ws.Range("A37"), $true)
$ws.Range("A57").Select()
